$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.762.70"
$ws.Range('E2').Value = "'  +0.42%  "
$ws.Range('D3').Value = "'1.918.24"
$ws.Range('E3').Value = "'  +1.61%  "
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = "'  -0.07%  "
$ws.Range('D5').Value = "'241.25"
$ws.Range('E5').Value = "'  -2.04%  "
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = "'  -0.02%  "
$ws.Range('D7').Value = "'0.4929"
$ws.Range('E7').Value = "'  +0.32%  "
$ws.Range('D8').Value = "'0.3028"
$ws.Range('E8').Value = "'  +2.83%  "
$ws.Range('D9').Value = "'0.06802"
$ws.Range('E9').Value = "'  +0.39%  "
$ws.Range('D10').Value = "'1.916.25"
$ws.Range('E10').Value = "'  +1.50%  "
$ws.Range('D11').Value = "'17.26"
$ws.Range('E11').Value = "'  +0.55%  "
$ws.Range('D12').Value = "'0.07338"
$ws.Range('E12').Value = "'  +1.37%  "
$ws.Range('E13').Value = "'  +3.09%  "
$ws.Range('D14').Value = "'89.16"
$ws.Range('E14').Value = "'  -2.16%  "
$ws.Range('D15').Value = "'0.6788"
$ws.Range('E15').Value = "'  +0.38%  "
$ws.Range('D16').Value = "'30.750.51"
$ws.Range('E16').Value = "'  +0.45%  "
$ws.Range('D17').Value = "'0.000008035"
$ws.Range('E17').Value = "'  +0.92%  "
$ws.Range('D18').Value = "'13.66"
$ws.Range('E18').Value = "'  +3.52%  "
$ws.Range('D19').Value = "'1.000"
$ws.Range('E19').Value = "'  -0.01%  "
$ws.Range('D20').Value = "'2.164.57"
$ws.Range('E20').Value = "'  +1.53%  "
$ws.Range('D21').Value = "'5.372"
$ws.Range('E21').Value = "'  +11.54%  "
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = "'  -0.06%  "
$ws.Range('D23').Value = "'202.74"
$ws.Range('E23').Value = "'  +10.54%  "
$ws.Range('D24').Value = "'6.340"
$ws.Range('E24').Value = "'  +4.62%  "
$ws.Range('D25').Value = "'9.739"
$ws.Range('E25').Value = "'  +4.00%  "
$ws.Range('D26').Value = "'161.44"
$ws.Range('E26').Value = "'  +3.88%  "
$ws.Range('D27').Value = "'18.87"
$ws.Range('E27').Value = "'  -0.75%  "
$ws.Range('E28').Value = "'  +3.63%  "
$ws.Range('D29').Value = "'1.454"
$ws.Range('E29').Value = "'  +3.99%  "
$ws.Range('D30').Value = "'4.391"
$ws.Range('E30').Value = "'  +1.69%  "
$ws.Range('E31').Value = "'  +1.93%  "
$ws.Range('D32').Value = "'4.121"
$ws.Range('E32').Value = "'  +2.88%  "
$ws.Range('D33').Value = "'0.05332"
$ws.Range('E33').Value = "'  +2.58%  "
$ws.Range('D34').Value = "'0.7511"
$ws.Range('E34').Value = "'  -0.33%  "
$ws.Range('D35').Value = "'1.129"
$ws.Range('E35').Value = "'  +1.48%  "
$ws.Range('D36').Value = "'2.698"
$ws.Range('E36').Value = "'  -1.77%  "
$ws.Range('D37').Value = "'0.01870"
$ws.Range('E37').Value = "'  +1.50%  "
$ws.Range('D38').Value = "'2.729"
$ws.Range('E38').Value = "'  +2.47%  "
$ws.Range('D39').Value = "'0.9314"
$ws.Range('E39').Value = "'  -0.70%  "
$ws.Range('D40').Value = "'2.094"
$ws.Range('E40').Value = "'  -2.48%  "
$ws.Range('D41').Value = "'0.4517"
$ws.Range('E41').Value = "'  +2.23%  "
$ws.Range('D42').Value = "'72.95"
$ws.Range('E42').Value = "'  +26.01%  "
$ws.Range('D43').Value = "'108.00"
$ws.Range('E43').Value = "'  +2.33%  "
$ws.Range('D44').Value = "'5.985"
$ws.Range('E44').Value = "'  +4.09%  "
$ws.Range('B45').Value = "'PaxDollar"
$ws.Range('C45').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('D45').Value = "'1.003"
$ws.Range('E45').Value = "'  +0.24%  "
$ws.Range('B46').Value = "'Algorand"
$ws.Range('C46').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('D46').Value = "'0.1405"
$ws.Range('E46').Value = "'  +4.93%  "
$ws.Range('D47').Value = "'7.796"
$ws.Range('E47').Value = "'  +2.55%  "
$ws.Range('D48').Value = "'35.99"
$ws.Range('E48').Value = "'  +7.21%  "
$ws.Range('D49').Value = "'9.156"
$ws.Range('E49').Value = "'  +5.22%  "
$ws.Range('D50').Value = "'0.05962"
$ws.Range('E50').Value = "'  +2.03%  "
$ws.Range('D51').Value = "'0.4070"
$ws.Range('E51').Value = "'  +3.52%  "
